$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dealer")

# Row 8: SEIS732_Team_02_Corporate / Owner_Address / OADR_County
$ws.Range("G8").Value = "SEIS732_Team_02_Corporate"
$ws.Range("H8").Value = "Owner_Address"
$ws.Range("I8").Value = "OADR_County"

# Row 10: SEIS732_Team_02_Corporate / Owner_Address / OADR_Country
$ws.Range("G10").Value = "SEIS732_Team_02_Corporate"
$ws.Range("H10").Value = "Owner_Address"
$ws.Range("I10").Value = "OADR_Country"

# Update the active cell selection to I11 as shown in the diff
$ws.Range("I11").Select()
